$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.085.96'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.50%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.551.19'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.93%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("E5").Value = '  -0.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '287.42'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.16%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3821'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.59%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3288'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.66%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.58'
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.129'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.04%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07347'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.49%  '

$ws.Range("E12").Value = '  -0.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.12'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.19%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.778'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.37%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.726'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.568.60'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.38%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001068'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.80%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06633'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.82%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '85.59'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.31%  '

$ws.Range("E20").Value = '  -0.07%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.339'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.03'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.69'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.95%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.094.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.42%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.300'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.47%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.525'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '150.29'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.82%  '

$ws.Range("E28").Value = '  -2.49%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.929'
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.77'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.78%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.741.16'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.20%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.069'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.62%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.874'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.75%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.898'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.78%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08211'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.11%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.284'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.40%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06314'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.88%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02314'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.52%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.265'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.94%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2140'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.61%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.230'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.07%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.97'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.17%  '

$ws.Range("E43").Value = '  -0.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6008'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.75'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.10%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.729'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5822'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.66%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.966'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.45%  '

$ws.Range("E49").Value = '  -2.86%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.171'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.01%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07010'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.92%  '
